$d = $word.ActiveDocument

# 1. The title paragraph currently reads "D" + a "_GoBack" bookmark + "ate",
#    split across two runs. Collapse it into a single run with the text
#    "Date" (this also removes the bookmark that was sitting in between).
$d.Content.Find.Execute("Date", $true, $false, $false, $false, $false, $true, 1, $false, "Date", 2)

# 2. Re-create the "_GoBack" bookmark on the (now) empty paragraph that
#    immediately follows the "Agenda" heading.
$idx = 0
$agendaIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text.TrimEnd("`r", "`a")
    if ($t -eq "Agenda") {
        $agendaIdx = $idx
        break
    }
}
$target = $d.Paragraphs.Item($agendaIdx + 1)
$d.Bookmarks.Add("_GoBack", $target.Range)
